# New weekly "Kiwi" price-sheet refresh for Terminal Hortofruticola Agro Chillan.
# Two new report rows (dated 2023-07-17 / serial 45124) are prepended to the
# existing "Especial"/"Primera" pair, pushing every later record down by two
# rows; the two oldest records that fall off the bottom (old rows 331-332) are
# re-appended as the new rows 333-334 so no historical data is lost.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rows 286-287: keep all data, only the report date moves forward ---
$ws.Range("D286").Value = 45124
$ws.Range("D287").Value = 45124

# --- Rows 288-332: each row now holds what used to be two rows earlier ---
# row 288 <- old row 286
$ws.Range("D288").Value = 45070
$ws.Range("L288").Value = 'Especial'
$ws.Range("N288").Value = 12000
$ws.Range("O288").Value = 12000
$ws.Range("P288").Value = 12000
$ws.Range("R288").Value = 'Región de O''Higgins'
$ws.Range("S288").Value = 667
# row 289 <- old row 287
$ws.Range("D289").Value = 45070
$ws.Range("L289").Value = 'Primera'
$ws.Range("M289").Value = 80
$ws.Range("N289").Value = 10000
$ws.Range("O289").Value = 10000
$ws.Range("P289").Value = 10000
$ws.Range("R289").Value = 'Región de O''Higgins'
$ws.Range("S289").Value = 556
# row 290 <- old row 288
$ws.Range("D290").Value = 44790
$ws.Range("N290").Value = 7500
$ws.Range("O290").Value = 8000
$ws.Range("P290").Value = 7750
$ws.Range("R290").Value = 'Provincia de Curicó'
$ws.Range("S290").Value = 431
# row 291 <- old row 289
$ws.Range("D291").Value = 44790
$ws.Range("L291").Value = 'Segunda'
$ws.Range("M291").Value = 120
$ws.Range("N291").Value = 6500
$ws.Range("O291").Value = 7000
$ws.Range("P291").Value = 6750
$ws.Range("S291").Value = 375
# row 292 <- old row 290
$ws.Range("D292").Value = 45040
$ws.Range("L292").Value = 'Primera'
$ws.Range("M292").Value = 80
$ws.Range("N292").Value = 13000
$ws.Range("O292").Value = 13000
$ws.Range("P292").Value = 13000
$ws.Range("R292").Value = 'Región de O''Higgins'
$ws.Range("S292").Value = 722
# row 293 <- old row 291
$ws.Range("D293").Value = 44692
$ws.Range("M293").Value = 160
$ws.Range("O293").Value = 11000
$ws.Range("P293").Value = 10500
$ws.Range("R293").Value = 'Provincia de Curicó'
$ws.Range("S293").Value = 583
# row 294 <- old row 292
$ws.Range("D294").Value = 44692
$ws.Range("M294").Value = 60
$ws.Range("O294").Value = 9000
$ws.Range("P294").Value = 8500
$ws.Range("R294").Value = 'Provincia de Curicó'
$ws.Range("S294").Value = 472
# row 295 <- old row 293
$ws.Range("D295").Value = 45111
$ws.Range("M295").Value = 60
$ws.Range("N295").Value = 10000
$ws.Range("O295").Value = 10000
$ws.Range("P295").Value = 10000
$ws.Range("R295").Value = 'Región de O''Higgins'
$ws.Range("S295").Value = 556
# row 296 <- old row 294
$ws.Range("D296").Value = 45111
$ws.Range("M296").Value = 40
$ws.Range("N296").Value = 8000
$ws.Range("O296").Value = 8000
$ws.Range("P296").Value = 8000
$ws.Range("R296").Value = 'Región de O''Higgins'
$ws.Range("S296").Value = 444
# row 297 <- old row 295
$ws.Range("D297").Value = 44417
$ws.Range("M297").Value = 160
$ws.Range("N297").Value = 12500
$ws.Range("O297").Value = 13000
$ws.Range("P297").Value = 12750
$ws.Range("Q297").Value = '$/bandeja 18 kilos'
$ws.Range("S297").Value = 708
$ws.Range("T297").Value = 18
# row 298 <- old row 296
$ws.Range("D298").Value = 44417
$ws.Range("L298").Value = 'Segunda'
$ws.Range("M298").Value = 80
$ws.Range("N298").Value = 11000
$ws.Range("O298").Value = 11500
$ws.Range("P298").Value = 11250
$ws.Range("R298").Value = 'Provincia de Curicó'
$ws.Range("S298").Value = 625
# row 299 <- old row 297
$ws.Range("D299").Value = 44323
$ws.Range("M299").Value = 120
$ws.Range("O299").Value = 11000
$ws.Range("P299").Value = 10500
$ws.Range("Q299").Value = '$/bandeja 10 kilos'
$ws.Range("R299").Value = 'Provincia de Curicó'
$ws.Range("S299").Value = 1050
$ws.Range("T299").Value = 10
# row 300 <- old row 298
$ws.Range("L300").Value = 'Especial'
$ws.Range("M300").Value = 60
$ws.Range("N300").Value = 12000
$ws.Range("O300").Value = 12000
$ws.Range("P300").Value = 12000
$ws.Range("S300").Value = 667
# row 301 <- old row 299
$ws.Range("D301").Value = 45106
$ws.Range("M301").Value = 40
$ws.Range("N301").Value = 10000
$ws.Range("O301").Value = 10000
$ws.Range("P301").Value = 10000
$ws.Range("R301").Value = 'Región de O''Higgins'
$ws.Range("S301").Value = 556
# row 302 <- old row 300
$ws.Range("D302").Value = 45106
$ws.Range("M302").Value = 40
$ws.Range("N302").Value = 8000
$ws.Range("O302").Value = 8000
$ws.Range("P302").Value = 8000
$ws.Range("R302").Value = 'Región de O''Higgins'
$ws.Range("S302").Value = 444
# row 303 <- old row 301
$ws.Range("D303").Value = 44358
$ws.Range("N303").Value = 10500
$ws.Range("P303").Value = 10750
$ws.Range("S303").Value = 597
# row 304 <- old row 302
$ws.Range("D304").Value = 44358
$ws.Range("L304").Value = 'Segunda'
$ws.Range("M304").Value = 120
$ws.Range("N304").Value = 8500
$ws.Range("O304").Value = 9000
$ws.Range("P304").Value = 8750
$ws.Range("S304").Value = 486
# row 305 <- old row 303
$ws.Range("D305").Value = 44391
$ws.Range("L305").Value = 'Primera'
$ws.Range("N305").Value = 10000
$ws.Range("O305").Value = 11000
$ws.Range("P305").Value = 10500
$ws.Range("S305").Value = 583
# row 306 <- old row 304
$ws.Range("D306").Value = 44420
$ws.Range("L306").Value = 'Primera'
$ws.Range("M306").Value = 200
$ws.Range("O306").Value = 12500
$ws.Range("P306").Value = 12250
$ws.Range("R306").Value = 'Provincia de Curicó'
$ws.Range("S306").Value = 681
# row 307 <- old row 305
$ws.Range("D307").Value = 44420
$ws.Range("L307").Value = 'Segunda'
$ws.Range("M307").Value = 120
$ws.Range("N307").Value = 11000
$ws.Range("O307").Value = 11500
$ws.Range("P307").Value = 11250
$ws.Range("R307").Value = 'Provincia de Curicó'
$ws.Range("S307").Value = 625
# row 308 <- old row 306
$ws.Range("D308").Value = 45075
$ws.Range("L308").Value = 'Especial'
$ws.Range("M308").Value = 40
$ws.Range("N308").Value = 12000
$ws.Range("O308").Value = 12000
$ws.Range("P308").Value = 12000
$ws.Range("Q308").Value = '$/bandeja 18 kilos'
$ws.Range("R308").Value = 'Región de O''Higgins'
$ws.Range("S308").Value = 667
$ws.Range("T308").Value = 18
# row 309 <- old row 307
$ws.Range("D309").Value = 45075
$ws.Range("L309").Value = 'Primera'
$ws.Range("M309").Value = 30
$ws.Range("N309").Value = 10000
$ws.Range("O309").Value = 10000
$ws.Range("P309").Value = 10000
$ws.Range("Q309").Value = '$/bandeja 18 kilos'
$ws.Range("R309").Value = 'Región de O''Higgins'
$ws.Range("S309").Value = 556
$ws.Range("T309").Value = 18
# row 310 <- old row 308
$ws.Range("D310").Value = 44364
$ws.Range("L310").Value = 'Primera'
$ws.Range("M310").Value = 120
$ws.Range("N310").Value = 8000
$ws.Range("O310").Value = 8500
$ws.Range("P310").Value = 8250
$ws.Range("Q310").Value = '$/bandeja 10 kilos'
$ws.Range("R310").Value = 'Provincia de Curicó'
$ws.Range("S310").Value = 825
$ws.Range("T310").Value = 10
# row 311 <- old row 309
$ws.Range("D311").Value = 44364
$ws.Range("L311").Value = 'Segunda'
$ws.Range("M311").Value = 120
$ws.Range("N311").Value = 7000
$ws.Range("O311").Value = 7500
$ws.Range("P311").Value = 7250
$ws.Range("Q311").Value = '$/bandeja 10 kilos'
$ws.Range("R311").Value = 'Provincia de Curicó'
$ws.Range("S311").Value = 725
$ws.Range("T311").Value = 10
# row 312 <- old row 310
$ws.Range("L312").Value = 'Especial'
$ws.Range("M312").Value = 60
$ws.Range("N312").Value = 11000
$ws.Range("O312").Value = 11000
$ws.Range("P312").Value = 11000
$ws.Range("S312").Value = 611
# row 313 <- old row 311
$ws.Range("M313").Value = 50
# row 314 <- old row 312
$ws.Range("M314").Value = 30
$ws.Range("N314").Value = 9000
$ws.Range("O314").Value = 9000
$ws.Range("P314").Value = 9000
$ws.Range("S314").Value = 500
# row 315 <- old row 313
$ws.Range("D315").Value = 45112
$ws.Range("M315").Value = 60
$ws.Range("O315").Value = 10000
$ws.Range("P315").Value = 10000
$ws.Range("R315").Value = 'Región de O''Higgins'
$ws.Range("S315").Value = 556
# row 316 <- old row 314
$ws.Range("D316").Value = 45112
$ws.Range("M316").Value = 40
$ws.Range("O316").Value = 8000
$ws.Range("P316").Value = 8000
$ws.Range("R316").Value = 'Región de O''Higgins'
$ws.Range("S316").Value = 444
# row 317 <- old row 315
$ws.Range("D317").Value = 44348
$ws.Range("N317").Value = 10000
$ws.Range("O317").Value = 11000
$ws.Range("P317").Value = 10500
$ws.Range("S317").Value = 583
# row 318 <- old row 316
$ws.Range("D318").Value = 44348
$ws.Range("M318").Value = 120
$ws.Range("N318").Value = 8000
$ws.Range("O318").Value = 9000
$ws.Range("P318").Value = 8500
$ws.Range("S318").Value = 472
# row 319 <- old row 317
$ws.Range("D319").Value = 44749
$ws.Range("L319").Value = 'Primera'
$ws.Range("M319").Value = 120
$ws.Range("N319").Value = 6500
$ws.Range("O319").Value = 7000
$ws.Range("P319").Value = 6750
$ws.Range("R319").Value = 'Provincia de Curicó'
$ws.Range("S319").Value = 375
# row 320 <- old row 318
$ws.Range("D320").Value = 44749
$ws.Range("L320").Value = 'Segunda'
$ws.Range("M320").Value = 60
$ws.Range("N320").Value = 5500
$ws.Range("O320").Value = 5500
$ws.Range("P320").Value = 5500
$ws.Range("R320").Value = 'Provincia de Curicó'
$ws.Range("S320").Value = 306
# row 321 <- old row 319
$ws.Range("L321").Value = 'Especial'
$ws.Range("M321").Value = 80
$ws.Range("N321").Value = 12000
$ws.Range("O321").Value = 12000
$ws.Range("P321").Value = 12000
$ws.Range("S321").Value = 667
# row 322 <- old row 320
$ws.Range("D322").Value = 45089
$ws.Range("M322").Value = 80
$ws.Range("N322").Value = 10000
$ws.Range("O322").Value = 10000
$ws.Range("P322").Value = 10000
$ws.Range("R322").Value = 'Región de O''Higgins'
$ws.Range("S322").Value = 556
# row 323 <- old row 321
$ws.Range("D323").Value = 45089
$ws.Range("N323").Value = 8000
$ws.Range("O323").Value = 8000
$ws.Range("P323").Value = 8000
$ws.Range("R323").Value = 'Región de O''Higgins'
$ws.Range("S323").Value = 444
# row 324 <- old row 322
$ws.Range("D324").Value = 44763
$ws.Range("L324").Value = 'Primera'
$ws.Range("M324").Value = 120
$ws.Range("N324").Value = 5500
$ws.Range("O324").Value = 6000
$ws.Range("P324").Value = 5750
$ws.Range("R324").Value = 'Provincia de Curicó'
$ws.Range("S324").Value = 319
# row 325 <- old row 323
$ws.Range("D325").Value = 44763
$ws.Range("L325").Value = 'Segunda'
$ws.Range("M325").Value = 60
$ws.Range("N325").Value = 5000
$ws.Range("O325").Value = 5000
$ws.Range("P325").Value = 5000
$ws.Range("R325").Value = 'Provincia de Curicó'
$ws.Range("S325").Value = 278
# row 326 <- old row 324
$ws.Range("D326").Value = 45099
$ws.Range("L326").Value = 'Especial'
$ws.Range("N326").Value = 12000
$ws.Range("O326").Value = 12000
$ws.Range("P326").Value = 12000
$ws.Range("S326").Value = 667
# row 327 <- old row 325
$ws.Range("D327").Value = 45099
$ws.Range("L327").Value = 'Primera'
$ws.Range("M327").Value = 100
$ws.Range("N327").Value = 10000
$ws.Range("O327").Value = 10000
$ws.Range("P327").Value = 10000
$ws.Range("S327").Value = 556
# row 328 <- old row 326
$ws.Range("D328").Value = 45121
$ws.Range("M328").Value = 80
$ws.Range("N328").Value = 10000
$ws.Range("O328").Value = 10000
$ws.Range("P328").Value = 10000
$ws.Range("R328").Value = 'Región de O''Higgins'
$ws.Range("S328").Value = 556
# row 329 <- old row 327
$ws.Range("D329").Value = 45121
$ws.Range("M329").Value = 50
$ws.Range("N329").Value = 8000
$ws.Range("O329").Value = 8000
$ws.Range("P329").Value = 8000
$ws.Range("R329").Value = 'Región de O''Higgins'
$ws.Range("S329").Value = 444
# row 330 <- old row 328
$ws.Range("D330").Value = 44777
$ws.Range("L330").Value = 'Primera'
$ws.Range("M330").Value = 120
$ws.Range("N330").Value = 6500
$ws.Range("O330").Value = 7000
$ws.Range("P330").Value = 6750
$ws.Range("R330").Value = 'Provincia de Curicó'
$ws.Range("S330").Value = 375
# row 331 <- old row 329
$ws.Range("D331").Value = 44777
$ws.Range("L331").Value = 'Segunda'
$ws.Range("M331").Value = 120
$ws.Range("N331").Value = 5000
$ws.Range("O331").Value = 5500
$ws.Range("P331").Value = 5250
$ws.Range("R331").Value = 'Provincia de Curicó'
$ws.Range("S331").Value = 292
# row 332 <- old row 330
$ws.Range("L332").Value = 'Especial'
$ws.Range("M332").Value = 40
$ws.Range("N332").Value = 12000
$ws.Range("O332").Value = 12000
$ws.Range("P332").Value = 12000
$ws.Range("S332").Value = 667

# --- Rows 333-334: brand-new rows carrying the displaced old rows 331-332 ---
# row 333 <- old row 331
$ws.Range("A333").Value = 7
$ws.Range("B333").Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Range("C333").Value = 'Ñuble'
$ws.Range("D333").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("D333").Value = 45072
$ws.Range("E333").Value = 16
$ws.Range("F333").Value = 'Fruta'
$ws.Range("G333").Value = 100101
$ws.Range("H333").Value = 'Berries'
$ws.Range("I333").Value = 100101007
$ws.Range("J333").Value = 'Kiwi'
$ws.Range("K333").Value = 'Hayward'
$ws.Range("L333").Value = 'Primera'
$ws.Range("M333").Value = 60
$ws.Range("N333").Value = 10000
$ws.Range("O333").Value = 10000
$ws.Range("P333").Value = 10000
$ws.Range("Q333").Value = '$/bandeja 18 kilos'
$ws.Range("R333").Value = 'Región de O''Higgins'
$ws.Range("S333").Value = 556
$ws.Range("T333").Value = 18
# row 334 <- old row 332
$ws.Range("A334").Value = 7
$ws.Range("B334").Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Range("C334").Value = 'Ñuble'
$ws.Range("D334").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("D334").Value = 45072
$ws.Range("E334").Value = 16
$ws.Range("F334").Value = 'Fruta'
$ws.Range("G334").Value = 100101
$ws.Range("H334").Value = 'Berries'
$ws.Range("I334").Value = 100101007
$ws.Range("J334").Value = 'Kiwi'
$ws.Range("K334").Value = 'Hayward'
$ws.Range("L334").Value = 'Segunda'
$ws.Range("M334").Value = 30
$ws.Range("N334").Value = 8000
$ws.Range("O334").Value = 8000
$ws.Range("P334").Value = 8000
$ws.Range("Q334").Value = '$/bandeja 18 kilos'
$ws.Range("R334").Value = 'Región de O''Higgins'
$ws.Range("S334").Value = 444
$ws.Range("T334").Value = 18
